$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.081251
$ws.Range("H2").Value = 9.243753
$ws.Range("I2").Value = 0.6013642694204734
$ws.Range("J2").Value = 0.6013642694204734
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 15.35884066666667
$ws.Range("N2").Value = 46.076522
$ws.Range("O2").Value = 0.1012042817263867
$ws.Range("P2").Value = 0.1012042817263867
$ws.Range("Q2").Value = 47.32444316300733
$ws.Range("R2").Value = 425.919988467066
$ws.Range("S2").Value = 0.06086063894261227
$ws.Range("T2").Value = 0.06086063894261228
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.081251
$ws.Range("H3").Value = 9.243753
$ws.Range("I3").Value = 0.6013642694204734
$ws.Range("J3").Value = 0.6013642694204734
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 50.59256466666667
$ws.Range("N3").Value = 151.777694
$ws.Range("O3").Value = 0.3333704853712116
$ws.Range("P3").Value = 0.3333704853712116
$ws.Range("Q3").Value = 155.8883904717313
$ws.Range("R3").Value = 1402.995514245582
$ws.Range("S3").Value = 0.2004770983816072
$ws.Range("T3").Value = 0.2004770983816072
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.081251
$ws.Range("H4").Value = 9.243753
$ws.Range("I4").Value = 0.6013642694204734
$ws.Range("J4").Value = 0.6013642694204734
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 60.37715666666667
$ws.Range("N4").Value = 181.13147
$ws.Range("O4").Value = 0.397844271305776
$ws.Range("P4").Value = 0.397844271305776
$ws.Range("Q4").Value = 186.0371743563233
$ws.Range("R4").Value = 1674.33456920691
$ws.Range("S4").Value = 0.2392493295569186
$ws.Range("T4").Value = 0.2392493295569186
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.081251
$ws.Range("H5").Value = 9.243753
$ws.Range("I5").Value = 0.6013642694204734
$ws.Range("J5").Value = 0.6013642694204734
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.43221733333333
$ws.Range("N5").Value = 76.296652
$ws.Range("O5").Value = 0.1675809615966257
$ws.Range("P5").Value = 0.1675809615966258
$ws.Range("Q5").Value = 78.36304509055066
$ws.Range("R5").Value = 705.267405814956
$ws.Range("S5").Value = 0.1007772025393353
$ws.Range("T5").Value = 0.1007772025393353
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.6679959999999999
$ws.Range("H6").Value = 2.003988
$ws.Range("I6").Value = 0.1303720230892577
$ws.Range("J6").Value = 0.1303720230892577
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 15.35884066666667
$ws.Range("N6").Value = 46.076522
$ws.Range("O6").Value = 0.1012042817263867
$ws.Range("P6").Value = 0.1012042817263867
$ws.Range("Q6").Value = 10.25964412997067
$ws.Range("R6").Value = 92.33679716973597
$ws.Range("S6").Value = 0.01319420695396423
$ws.Range("T6").Value = 0.01319420695396423
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.6679959999999999
$ws.Range("H7").Value = 2.003988
$ws.Range("I7").Value = 0.1303720230892577
$ws.Range("J7").Value = 0.1303720230892577
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 50.59256466666667
$ws.Range("N7").Value = 151.777694
$ws.Range("O7").Value = 0.3333704853712116
$ws.Range("P7").Value = 0.3333704853712116
$ws.Range("Q7").Value = 33.79563082707466
$ws.Range("R7").Value = 304.1606774436719
$ws.Range("S7").Value = 0.04346218461609266
$ws.Range("T7").Value = 0.04346218461609266
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.6679959999999999
$ws.Range("H8").Value = 2.003988
$ws.Range("I8").Value = 0.1303720230892577
$ws.Range("J8").Value = 0.1303720230892577
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 60.37715666666667
$ws.Range("N8").Value = 181.13147
$ws.Range("O8").Value = 0.397844271305776
$ws.Range("P8").Value = 0.397844271305776
$ws.Range("Q8").Value = 40.33169914470667
$ws.Range("R8").Value = 362.9852923023599
$ws.Range("S8").Value = 0.05186776252460555
$ws.Range("T8").Value = 0.05186776252460555
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.6679959999999999
$ws.Range("H9").Value = 2.003988
$ws.Range("I9").Value = 0.1303720230892577
$ws.Range("J9").Value = 0.1303720230892577
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 25.43221733333333
$ws.Range("N9").Value = 76.296652
$ws.Range("O9").Value = 0.1675809615966257
$ws.Range("P9").Value = 0.1675809615966258
$ws.Range("Q9").Value = 16.98861944979733
$ws.Range("R9").Value = 152.897575048176
$ws.Range("S9").Value = 0.02184786899459531
$ws.Range("T9").Value = 0.02184786899459531
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.374521
$ws.Range("H10").Value = 4.123563
$ws.Range("I10").Value = 0.2682637074902688
$ws.Range("J10").Value = 0.2682637074902689
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 15.35884066666667
$ws.Range("N10").Value = 46.076522
$ws.Range("O10").Value = 0.1012042817263867
$ws.Range("P10").Value = 0.1012042817263867
$ws.Range("Q10").Value = 21.11104903198733
$ws.Range("R10").Value = 189.999441287886
$ws.Range("S10").Value = 0.02714943582981015
$ws.Range("T10").Value = 0.02714943582981016
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.374521
$ws.Range("H11").Value = 4.123563
$ws.Range("I11").Value = 0.2682637074902688
$ws.Range("J11").Value = 0.2682637074902689
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 50.59256466666667
$ws.Range("N11").Value = 151.777694
$ws.Range("O11").Value = 0.3333704853712116
$ws.Range("P11").Value = 0.3333704853712116
$ws.Range("Q11").Value = 69.54054257819134
$ws.Range("R11").Value = 625.864883203722
$ws.Range("S11").Value = 0.08943120237351164
$ws.Range("T11").Value = 0.08943120237351167
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.374521
$ws.Range("H12").Value = 4.123563
$ws.Range("I12").Value = 0.2682637074902688
$ws.Range("J12").Value = 0.2682637074902689
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 60.37715666666667
$ws.Range("N12").Value = 181.13147
$ws.Range("O12").Value = 0.397844271305776
$ws.Range("P12").Value = 0.397844271305776
$ws.Range("Q12").Value = 82.98966975862334
$ws.Range("R12").Value = 746.90702782761
$ws.Range("S12").Value = 0.1067271792242518
$ws.Range("T12").Value = 0.1067271792242519
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.374521
$ws.Range("H13").Value = 4.123563
$ws.Range("I13").Value = 0.2682637074902688
$ws.Range("J13").Value = 0.2682637074902689
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 25.43221733333333
$ws.Range("N13").Value = 76.296652
$ws.Range("O13").Value = 0.1675809615966257
$ws.Range("P13").Value = 0.1675809615966258
$ws.Range("Q13").Value = 34.95711680123066
$ws.Range("R13").Value = 314.6140512110759
$ws.Range("S13").Value = 0.04495589006269519
$ws.Range("T13").Value = 0.0449558900626952